$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.161853194236755
$ws.Range("B1").Value = 2.372714042663574
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.392551898956299
$ws.Range("E1").Value = 1.218700766563416
